$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8063246238521696, 1.1157283852016282, 0.9274878271935929, 2.8006964598969546)"
$ws.Range("C2").Value = "JSU(-0.868993742167784, 1.2036190664269486, 4.464934356795796, 5.903238932687726)"
$ws.Range("D2").Value = "NIG(0.9289038427542619, 0.7364796381811188, 0.754590683413744, 3.1094613782342373)"
$ws.Range("E2").Value = "NIG(1.911271607600541, 1.5553807283435634, 2.2304123877282884, 6.904349453559232)"
